$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New goods label order (ties in the count column were re-shuffled by a re-run
# of the pandas export referenced in the commit message). Counts in column B are
# untouched; only the word in column A is rewritten for each row.
$ws.Cells.Item(2, 1).Value = "хлеб"
$ws.Cells.Item(3, 1).Value = "вино"
$ws.Cells.Item(4, 1).Value = "скот"
$ws.Cells.Item(5, 1).Value = "холст"
$ws.Cells.Item(6, 1).Value = "кожа"
$ws.Cells.Item(7, 1).Value = "мед"
$ws.Cells.Item(8, 1).Value = "пиво"
$ws.Cells.Item(9, 1).Value = "сукно"
$ws.Cells.Item(10, 1).Value = "овчина"
$ws.Cells.Item(11, 1).Value = "лошадь"
$ws.Cells.Item(12, 1).Value = "воск"
$ws.Cells.Item(13, 1).Value = "масло"
$ws.Cells.Item(14, 1).Value = "сало"
$ws.Cells.Item(15, 1).Value = "железо"
$ws.Cells.Item(16, 1).Value = "Крымскую соль"
$ws.Cells.Item(17, 1).Value = "полотно"
$ws.Cells.Item(18, 1).Value = "колеса"
$ws.Cells.Item(19, 1).Value = "сено"
$ws.Cells.Item(20, 1).Value = "парча"
$ws.Cells.Item(21, 1).Value = "говядина"
$ws.Cells.Item(22, 1).Value = "табак"
$ws.Cells.Item(23, 1).Value = "позумент"
$ws.Cells.Item(24, 1).Value = "выбойка"
$ws.Cells.Item(25, 1).Value = "шелк"
$ws.Cells.Item(26, 1).Value = "чулок"
$ws.Cells.Item(27, 1).Value = "сахар"
$ws.Cells.Item(28, 1).Value = "лес"
$ws.Cells.Item(29, 1).Value = "лыко"
$ws.Cells.Item(30, 1).Value = "китайка"
$ws.Cells.Item(31, 1).Value = "сапог"
$ws.Cells.Item(32, 1).Value = "коса"
$ws.Cells.Item(33, 1).Value = "сани"
$ws.Cells.Item(34, 1).Value = "ладан"
$ws.Cells.Item(35, 1).Value = "гвоздь"
$ws.Cells.Item(36, 1).Value = "горшок"
$ws.Cells.Item(37, 1).Value = "ром"
$ws.Cells.Item(38, 1).Value = "конь"
$ws.Cells.Item(39, 1).Value = "обод"
$ws.Cells.Item(40, 1).Value = "рогожа"
$ws.Cells.Item(41, 1).Value = "платок"
$ws.Cells.Item(42, 1).Value = "овца"
$ws.Cells.Item(43, 1).Value = "замок"
$ws.Cells.Item(44, 1).Value = "веревка"
$ws.Cells.Item(45, 1).Value = "сосуд"
$ws.Cells.Item(46, 1).Value = "нитка"
$ws.Cells.Item(47, 1).Value = "гумми"
$ws.Cells.Item(48, 1).Value = "котел"
$ws.Cells.Item(49, 1).Value = "роза"
$ws.Cells.Item(50, 1).Value = "покроми"
$ws.Cells.Item(51, 1).Value = "брусья"
$ws.Cells.Item(52, 1).Value = "дуга"
$ws.Cells.Item(53, 1).Value = "бечева"
$ws.Cells.Item(54, 1).Value = "сковорода"
$ws.Cells.Item(55, 1).Value = "хомут"
$ws.Cells.Item(56, 1).Value = "скотский кожа"

